# items.xlsx — update Loot-Box related item descriptions, rename a couple
# of materials' purposes, add two new "Loot Box" items (rows 14 & 15),
# and move the active selection/view to the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Description" (column B) text for rows 7-13 ---------
$ws.Range("B7").Value  = "Dropped by Golem. Material exchanged to obtain Common Loot Boxes"
$ws.Range("B8").Value  = "Dropped by Golem. Material exchanged to obtain Premium Loot Boxes"
$ws.Range("B9").Value  = "Dropped by Panther. Used to restore health"
$ws.Range("B10").Value = "Dropped by Panther. Material exchanged to obtain Common Loot Boxes"
$ws.Range("B11").Value = "Dropped by Panther. Material exchanged to obtain Premium Loot Boxes"
$ws.Range("B12").Value = "Dropped by Treant. Material exchanged to obtain Common Loot Boxes"
$ws.Range("B13").Value = "Dropped by Treant. Material exchanged to obtain Premium Loot Boxes"

# --- Row heights now that the wrapped description text has changed -------
$ws.Rows.Item(7).RowHeight  = 72
$ws.Rows.Item(9).RowHeight  = 43.2
$ws.Rows.Item(10).RowHeight = 72
$ws.Rows.Item(12).RowHeight = 72

# --- New rows: Common / Premium Loot Box items ----------------------------
$ws.Range("A14").Value = "Common Loot Box"
$ws.Range("B14").Value = "Exchanged for using common materials. Can drop exclusive armor, weapons, and pets"
$ws.Range("C14").Value = "None"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "Loot Box"
$ws.Rows.Item(14).RowHeight = 86.4

$ws.Range("A15").Value = "Premium Loot Box"
$ws.Range("B15").Value = "Exchanged for using premium materials. Can drop exclusive armor, weapons, and pets"
$ws.Range("C15").Value = "None"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "Loot Box"
$ws.Rows.Item(15).RowHeight = 86.4

# --- Move the view / selection down to the newly added rows --------------
$ws.Range("H14").Select() | Out-Null
